$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add three new rows of data (17, 18, 19) with Name in column A and notes
# in column C. The fill order below matches the order in which new
# strings were appended to the shared string table.
$ws.Range("A17").Value = "eucharia"
$ws.Range("A18").Value = "Aemelius Magnus Arborius"
$ws.Range("C17").Value = "Lamaire PLM v2"
$ws.Range("C18").Value = "Lamaire PLM v2"
$ws.Range("A19").Value = "Gaius Cassius Parmensis"
$ws.Range("C19").Value = "Lemaire PLM v"

# Update the view to reflect the new selection/scroll position
$ws.Range("C19").Select()
$excel.ActiveWindow.ScrollRow = 4
